$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.523.27'
$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("D3").Value = '1.951.84'
$ws.Range("E3").Value = '  +0.70%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.625'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.79'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.95%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.377'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0786'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.85%  '
$ws.Range("E11").Value = '  +0.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.12'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.71%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.836'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.32%  '
$ws.Range("D14").Value = '2.240.29'
$ws.Range("E14").Value = '  +0.83%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.43'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.25'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.56%  '
$ws.Range("D17").Value = '1.959.13'
$ws.Range("E17").Value = '  +1.25%  '
$ws.Range("D18").Value = '36.437.95'
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.11'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").Value = '0.0₃0851'
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '229.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.86%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.36'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.68%  '
$ws.Range("E26").Value = '  +7.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.14'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.39%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '160.20'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.22'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.94%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.30'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +20.40%  '
$ws.Range("E31").Value = '  +2.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.75'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0610'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.44'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.76%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.25'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.19%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.39'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.78'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.43'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -10.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0962'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.23%  '
$ws.Range("E41").Value = '  +0.18%  '
$ws.Range("E42").Value = '  +2.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0209'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.33%  '
$ws.Range("D45").Value = '1.358.31'
$ws.Range("E45").Value = '  +2.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.45'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.47%  '
$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.15'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.43%  '
$ws.Range("E49").Value = '  +0.91%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '46.11'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.67%  '
$ws.Range("D51").Value = '2.135.92'
$ws.Range("E51").Value = '  +1.05%  '
